$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Tier" column (N) entirely -------------------------------
$ws.Range("N1:N3").EntireColumn.Delete()

# --- Column widths ---------------------------------------------------------
$ws.Range("C1").EntireColumn.ColumnWidth = 35.5
$ws.Range("E1").EntireColumn.ColumnWidth = 11.83
$ws.Range("H1:J1").EntireColumn.ColumnWidth = 15.166666666666666

# --- Header row (row 1) -----------------------------------------------------
$ws.Range("A1").Value = "Order Received Data and Time"
$ws.Range("B1").Value = "OrderID"
$ws.Range("C1").Value = "Emp ID-Order Assigned"
$ws.Range("D1").Value = "Assignee_QA"
$ws.Range("E1").Value = "Client"
$ws.Range("F1").Value = "Typist"
$ws.Range("G1").Value = "Typist QC"
$ws.Range("H1").Value = "Product Name"
$ws.Range("I1").Value = "Process"
$ws.Range("J1").Value = "Lob"
$ws.Range("K1").Value = "State"
$ws.Range("L1").Value = "County"
$ws.Range("M1").Value = "Status"

# --- Row 2 -------------------------------------------------------------------
$ws.Range("A2").Value = 45436.041666666664
$ws.Range("B2").Value = 121321783
$ws.Range("C2").Value = "SIPL0005"
$ws.Range("D2").Value = "SIPL0004"
$ws.Range("E2").Value = "FINN TITLE"
$ws.Range("F2").Value = "SIPL0102"
$ws.Range("G2").Value = "SIPL0103"
$ws.Range("H2").Value = "Property Reports"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "Title"
$ws.Range("K2").Value = "FL"
$ws.Range("L2").Value = "Clay"
$ws.Range("M2").Value = "WIP"

# --- Row 3 -------------------------------------------------------------------
$ws.Range("A3").Value = 45439.083333333336
$ws.Range("B3").Value = 2193218321
$ws.Range("C3").Value = "SIPL0005"
$ws.Range("D3").Value = "SIPL0004"
$ws.Range("E3").Value = "FINN TITLE"
$ws.Range("F3").Value = "SIPL0102"
$ws.Range("G3").Value = "SIPL0103"
$ws.Range("H3").Value = "Foreclosure information Report"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Title"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "WIP"

# --- Fix cell formatting so it matches the pre-existing style pool --------
# The highlighted "Client" style (currently sitting on G2:G3) now belongs
# on column E only - move it there first, before it gets overwritten below.
$ws.Range("G2").Copy()
$ws.Range("E2:E3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Plain data-cell style (4-sided thin border, regular font) onto every other
# data cell that used to carry a distinctive highlight style at its old
# position, plus D2/D3 which used to carry a border missing the left edge.
$ws.Range("B2").Copy()
$ws.Range("C2:D3").PasteSpecial(-4122)
$ws.Range("F2:F3").PasteSpecial(-4122)
$ws.Range("G2:G3").PasteSpecial(-4122)
$ws.Range("H2:M3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Misc cosmetic selection ------------------------------------------------
$ws.Range("H12").Select()
